$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Ignore" conditional formatting rules that were on G2:H17
$cfRange = $ws.Range("G2:H17")
$cfRange.FormatConditions.Delete()

# Update the sheet view: scroll back to the top (drop the stale
# topLeftCell="A4") and move the active selection to K7
[void]$ws.Range("K7").Select()
